$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: keep only A1 ("psno"); drop the B1:D1 header cells entirely.
$ws.Range("B1:D1").ClearContents()

# Row 2: update A2 (keep it as text, not a number) and fill B2:U2.
$a2 = $ws.Range("A2")
$a2.NumberFormat = "@"
$a2.Value = "99004400.0"
$a2.Style = "Normal"

$ws.Range("B2").Value = "mordern family"
$ws.Range("C2").Value = "cyberpunk"
$ws.Range("D2").Value = "mirzapur"
$ws.Range("E2").Value = "anupama"
$ws.Range("F2").Value = "cyberpunk"
$ws.Range("G2").Value = "ac origins"
$ws.Range("H2").Value = "mordern family"
$ws.Range("I2").Value = "cyberpunk 2077"
$ws.Range("J2").Value = "ac valhalla"
$ws.Range("K2").Value = "resident evil"
$ws.Range("L2").Value = "cyberpunk"
$ws.Range("M2").Value = "ac valhalla"
$ws.Range("N2").Value = "writer's legacy"
$ws.Range("O2").Value = "two"
$ws.Range("P2").Value = "ac origins"
$ws.Range("Q2").Value = "resident evil5 "
$ws.Range("R2").Value = "euphoria"
$ws.Range("S2").Value = "mare easttown"
$ws.Range("T2").Value = "friends"
$ws.Range("U2").Value = "two"
